$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column (H1) — copy the formatting used by the other
# header cells (e.g. G1: bold font, border, centered) onto H1, then set
# its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Corresponding data value for the new column
$ws.Range("H2").Value = 1
